$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: link the local project with git ---
# A3 date stays 43480 (2019-01-15), unchanged
$ws.Range("C3").Value = "OLX website"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "link the local project with git"
$ws.Range("F3").Value = "10min"
$ws.Range("G3").Value = "5min"
$ws.Range("I3").Value = "Done"

# Action item cell wraps text, and row height grows to fit (matches row 2 style)
$ws.Range("E3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 30

# --- Row 4: basic smoke test ---
$ws.Range("A4").Value = 43481
$ws.Range("C4").Value = "OLX website"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "basic smoke test"
$ws.Range("F4").Value = "30min."

# --- Row 5: new project entry (OLX website, priority 1) ---
$ws.Range("A5").Value = 43481
$ws.Range("C5").Value = "OLX website"
$ws.Range("D5").Value = 1

# Update the active selection to F6 (matches the author's last edit position)
$ws.Range("F6").Select()
